$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: permanently split the run at a given absolute character position
# (without leaving any visible artifact) by adding a bookmark there and
# immediately deleting it again -- the run boundary that the bookmark forces
# survives the deletion.
# ---------------------------------------------------------------------------
function Split-RunAt($pos) {
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TEMP_SPLIT_MARKER", $bmRange) | Out-Null
    $d.Bookmarks("TEMP_SPLIT_MARKER").Delete()
}

# ---------------------------------------------------------------------------
# Step 1: delete the whole "Given the dataset ... radiotherapy?  (3 Marks)"
# paragraph (original question 3). Its question text is about to become the
# new question 2, and its "paragraph shell" is no longer needed.
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.Delete()

# After the deletion the paragraphs are:
#   3 -> "What percentage of patients received chemotherapy treatment? (2 Marks)"
#   4 -> "How many patients are in each tumor stage category? (2 Marks)"
#   5 -> "Create a new binary variable ... ER, PR, a[_GoBack]nd HER2 status. (3 Marks)"

# ---------------------------------------------------------------------------
# Step 2: question 1 becomes what used to be question 2's text.
# ---------------------------------------------------------------------------
$r1 = $d.Paragraphs(3).Range
$null = $r1.Find.Execute(
    "What percentage of patients received chemotherapy treatment?", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "How many patients are in each tumor stage category?", 2)

# Restore the run break right before " (2 Marks)" that the replace merged away.
$p1 = $d.Paragraphs(3).Range
$pos = $p1.Start + $p1.Text.IndexOf(" (2 Marks)")
Split-RunAt $pos

# ---------------------------------------------------------------------------
# Step 3: question 2 becomes what used to be question 3's text, and its
# marks run becomes "  (2 Marks)" (double space, matching the donor
# paragraph's formatting) with the relocated "_GoBack" bookmark sitting
# between the "2" and " Marks)".
# ---------------------------------------------------------------------------
$r2 = $d.Paragraphs(4).Range
$null = $r2.Find.Execute(
    "How many patients are in each tumor stage category? (2 Marks)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Given the dataset, how would you filter out patients who did not receive any chemotherapy or radiotherapy?  (2 Marks)", 2)

$p2 = $d.Paragraphs(4).Range
$p2Text = $p2.Text

# boundary between the question text and the "  (" run
$posA = $p2.Start + $p2Text.IndexOf("  (")
Split-RunAt $posA

# boundary between the "  (" run and the "2" run
$posB = $p2.Start + $p2Text.IndexOf("(2") + 1
Split-RunAt $posB

# boundary between the "2" run and the " Marks)" run -- this one keeps a
# real (relocated) "_GoBack" bookmark.
$posC = $p2.Start + $p2Text.IndexOf("2 Marks)") + 1
$bmRange = $d.Range($posC, $posC)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# Step 4: the old question 4 ("Create a new binary variable ...") loses the
# "_GoBack" bookmark that used to split "a" / "nd HER2 status." into two
# runs, and those two runs become one -- while the trailing " (3 Marks)"
# run stays separate.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$r3 = $d.Paragraphs(5).Range
$null = $r3.Find.Execute(
    "ER, PR, and HER2 status.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "ER, PR, and HER2 status.", 2)

$p3 = $d.Paragraphs(5).Range
$p3Text = $p3.Text
$posD = $p3.Start + $p3Text.IndexOf(" (3 Marks)")
Split-RunAt $posD

Write-Host "Done"
